$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 17.8699785029566
$arr[1,0] = 17.67362667991828
$arr[2,0] = 17.55474892445726
$arr[3,0] = 17.50677787488584
$arr[4,0] = 17.49884228981686
$arr[5,0] = 17.55409999194375
$arr[6,0] = 17.80195692414324
$arr[7,0] = 18.29915930166846
$arr[8,0] = 18.6682905735293
$arr[9,0] = 18.83642345784168
$arr[10,0] = 18.90006980480795
$arr[11,0] = 18.8863642277619
$arr[12,0] = 18.84166041847383
$arr[13,0] = 18.81427364389125
$arr[14,0] = 18.65730246900812
$arr[15,0] = 18.56102265905037
$arr[16,0] = 18.50566822623436
$arr[17,0] = 18.48693176206499
$arr[18,0] = 18.57126979520246
$arr[19,0] = 18.85479202047622
$arr[20,0] = 19.03993868822059
$arr[21,0] = 18.94115341075419
$arr[22,0] = 18.56663706676667
$arr[23,0] = 18.16376456175705
$ws.Range("B2:B25").Value2 = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 10.09996978357739
$arr[1,0] = 10.10525268607638
$arr[2,0] = 10.10893835428061
$arr[3,0] = 10.11055173908945
$arr[4,0] = 10.11082638061806
$arr[5,0] = 10.10895966128635
$arr[6,0] = 10.10169979129704
$arr[7,0] = 10.0909562411895
$arr[8,0] = 10.08517444156549
$arr[9,0] = 10.08299882290138
$arr[10,0] = 10.08224001248902
$arr[11,0] = 10.0824005478425
$arr[12,0] = 10.08293509323827
$arr[13,0] = 10.08327097968093
$arr[14,0] = 10.08532574322597
$arr[15,0] = 10.08670247210745
$arr[16,0] = 10.08753713360774
$arr[17,0] = 10.08782709673446
$arr[18,0] = 10.08655148925741
$arr[19,0] = 10.08277632139961
$arr[20,0] = 10.08068804003358
$arr[21,0] = 10.08176802013761
$arr[22,0] = 10.08661961415918
$arr[23,0] = 10.09349055217529
$ws.Range("D2:D25").Value2 = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 16.30835049313773
$arr[1,0] = 16.28214903295978
$arr[2,0] = 16.26657936079419
$arr[3,0] = 16.26036728846966
$arr[4,0] = 16.2593438624727
$arr[5,0] = 16.26649504217873
$arr[6,0] = 16.29921014527744
$arr[7,0] = 16.36737004670641
$arr[8,0] = 16.41974810623654
$arr[9,0] = 16.44404500175542
$arr[10,0] = 16.45331040568731
$arr[11,0] = 16.45131210391614
$arr[12,0] = 16.44480599208026
$arr[13,0] = 16.4408291502695
$arr[14,0] = 16.41816945758958
$arr[15,0] = 16.40438678571344
$arr[16,0] = 16.39650375685456
$arr[17,0] = 16.39384241096317
$arr[18,0] = 16.40584940012028
$arr[19,0] = 16.44671526206567
$arr[20,0] = 16.47379920564995
$arr[21,0] = 16.45931057421564
$arr[22,0] = 16.40518802482791
$arr[23,0] = 16.34851649820656
$ws.Range("E2:E25").Value2 = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 33.64749922674186
$arr[1,0] = 33.68332651649292
$arr[2,0] = 33.71405499980609
$arr[3,0] = 33.72876829738913
$arr[4,0] = 33.73134363110301
$arr[5,0] = 33.71424456269201
$arr[6,0] = 33.65803829586375
$arr[7,0] = 33.61723246380292
$arr[8,0] = 33.62970210189143
$arr[9,0] = 33.64459766265493
$arr[10,0] = 33.65156283916519
$arr[11,0] = 33.65000389423034
$arr[12,0] = 33.64514416008099
$arr[13,0] = 33.6423398566279
$arr[14,0] = 33.62891411470329
$arr[15,0] = 33.62303926921568
$arr[16,0] = 33.6205285106485
$arr[17,0] = 33.61982757422619
$arr[18,0] = 33.62357479957163
$arr[19,0] = 33.64653565354926
$arr[20,0] = 33.66926074133595
$arr[21,0] = 33.65642654015493
$arr[22,0] = 33.62332998620943
$arr[23,0] = 33.62082206283629
$ws.Range("F2:F25").Value2 = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 3.648505521904759
$arr[1,0] = 3.650984364246237
$arr[2,0] = 3.652588706645675
$arr[3,0] = 3.65326325755458
$arr[4,0] = 3.653376522563061
$arr[5,0] = 3.65259771969557
$arr[6,0] = 3.649343181261103
$arr[7,0] = 3.643611164810932
$arr[8,0] = 3.639791892213535
$arr[9,0] = 3.638138621189799
$arr[10,0] = 3.637524599531397
$arr[11,0] = 3.637656305756369
$arr[12,0] = 3.638087864390578
$arr[13,0] = 3.638353772140759
$arr[14,0] = 3.639901625038101
$arr[15,0] = 3.640872687196105
$arr[16,0] = 3.64143913936916
$arr[17,0] = 3.641632293015598
$arr[18,0] = 3.640768496434146
$arr[19,0] = 3.637960778993812
$arr[20,0] = 3.636195900910964
$arr[21,0] = 3.637131453358482
$arr[22,0] = 3.640815575606484
$arr[23,0] = 3.645092672455578
$ws.Range("G2:G25").Value2 = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 11.67327852520763
$arr[1,0] = 11.65999430185235
$arr[2,0] = 11.65216418490217
$arr[3,0] = 11.64905659476449
$arr[4,0] = 11.64854564308557
$arr[5,0] = 11.65212193627612
$arr[6,0] = 11.66863106767641
$arr[7,0] = 11.70354538728351
$arr[8,0] = 11.73067875305656
$arr[9,0] = 11.74332965842626
$arr[10,0] = 11.74816317990592
$arr[11,0] = 11.74712030860572
$arr[12,0] = 11.74372646329009
$arr[13,0] = 11.74165318081453
$arr[14,0] = 11.72985806435515
$arr[15,0] = 11.72270001479278
$arr[16,0] = 11.71861194868669
$arr[17,0] = 11.71723283108565
$arr[18,0] = 11.72345900414802
$arr[19,0] = 11.74472216465642
$arr[20,0] = 11.75886826487912
$arr[21,0] = 11.75129586124606
$arr[22,0] = 11.72311578026623
$arr[23,0] = 11.69383500309114
$ws.Range("J2:J25").Value2 = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 10.04245560319056
$arr[1,0] = 9.911965854353658
$arr[2,0] = 9.831918344966873
$arr[3,0] = 9.799349337398858
$arr[4,0] = 9.793945308545146
$arr[5,0] = 9.831478858519285
$arr[6,0] = 9.997466880474597
$arr[7,0] = 10.32213719911338
$arr[8,0] = 10.5582945153705
$arr[9,0] = 10.66482689374492
$arr[10,0] = 10.70500780837657
$arr[11,0] = 10.69636176444
$arr[12,0] = 10.66813599773738
$arr[13,0] = 10.65082506605877
$arr[14,0] = 10.55131160223703
$arr[15,0] = 10.49001081758735
$arr[16,0] = 10.45466965298221
$arr[17,0] = 10.44269053943915
$arr[18,0] = 10.49654516918142
$arr[19,0] = 10.67643119494214
$arr[20,0] = 10.79304715911242
$arr[21,0] = 10.73090412992331
$arr[22,0] = 10.49359129328546
$arr[23,0] = 9.799349337398858
$ws.Range("L2:L25").Value2 = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 16.57342728347234
$arr[1,0] = 16.45813053362408
$arr[2,0] = 16.38891123403293
$arr[3,0] = 16.36112130947402
$arr[4,0] = 16.35653265604983
$arr[5,0] = 16.38853473055421
$arr[6,0] = 16.53335971315126
$arr[7,0] = 16.82885875257273
$arr[8,0] = 17.05161919323529
$arr[9,0] = 17.1538920401442
$arr[10,0] = 17.19272979429907
$arr[11,0] = 17.18436093909625
$arr[12,0] = 17.15708520315374
$arr[13,0] = 17.14039149244845
$arr[14,0] = 17.04495210751932
$arr[15,0] = 16.98662517031603
$arr[16,0] = 16.95316689063242
$arr[17,0] = 16.94185470049202
$arr[18,0] = 16.99282506053729
$arr[19,0] = 17.16509398767059
$arr[20,0] = 17.27830654864444
$arr[21,0] = 17.21783418046154
$arr[22,0] = 16.99002185564999
$arr[23,0] = 16.74782548030715
$ws.Range("M2:M25").Value2 = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 19.80873072543924
$arr[1,0] = 19.85050309037126
$arr[2,0] = 19.87831702346665
$arr[3,0] = 19.89019687013348
$arr[4,0] = 19.89220247963404
$arr[5,0] = 19.87847502956037
$arr[6,0] = 19.82268492138283
$arr[7,0] = 19.73042371207206
$arr[8,0] = 19.6730373166148
$arr[9,0] = 19.64917748371305
$arr[10,0] = 19.64046441683709
$arr[11,0] = 19.64232661690075
$arr[12,0] = 19.6484542028063
$arr[13,0] = 19.65224945468733
$arr[14,0] = 19.67464173261773
$arr[15,0] = 19.68895325371729
$arr[16,0] = 19.69739626738398
$arr[17,0] = 19.70029125934208
$arr[18,0] = 19.68740789379817
$arr[19,0] = 19.64664564578007
$arr[20,0] = 19.62188244968735
$arr[21,0] = 19.63492751464148
$arr[22,0] = 19.68810588083548
$arr[23,0] = 19.75355302929979
$ws.Range("N2:N25").Value2 = $arr

$arr = New-Object 'object[,]' 24,1
$arr[0,0] = 25.24201608775375
$arr[1,0] = 25.24054338062222
$arr[2,0] = 25.24612494304928
$arr[3,0] = 25.25002996248794
$arr[4,0] = 25.25077683249939
$arr[5,0] = 25.24617100710361
$arr[6,0] = 25.24016222095304
$arr[7,0] = 25.27981644803055
$arr[8,0] = 25.34021113920844
$arr[9,0] = 25.3744304854352
$arr[10,0] = 25.38835281752919
$arr[11,0] = 25.38531161508752
$arr[12,0] = 25.37555659711882
$arr[13,0] = 25.36970674404557
$arr[14,0] = 25.3381100927222
$arr[15,0] = 25.32045050255566
$arr[16,0] = 25.31092851095949
$arr[17,0] = 25.30781380559018
$arr[18,0] = 25.32226467925475
$arr[19,0] = 25.37839576581235
$arr[20,0] = 25.42069733028641
$arr[21,0] = 25.39760847176684
$arr[22,0] = 25.32144252477723
$arr[23,0] = 25.26358775562582
$ws.Range("O2:O25").Value2 = $arr
